$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Content fixes for watch-list test cases (TestCase_E1, TestCase_E2 descriptions
# updated to reflect "ALL content" search results flow; TestCase_E3 result flips
# from SKIP to FAIL) ---

$ws1.Cells.Item(2, 3).Value = "Verify that user is able to add document to watchlist from ALL content search results page"
$ws1.Cells.Item(3, 3).Value = "Verify that user is able to add document to watchlist from document page once it is opened from ALL content set results"
$ws1.Cells.Item(4, 5).Value = "FAIL"

# --- Column width tweaks on the Test Cases sheet (C widened to fit the longer
# description text, E nudged slightly) ---
$ws1.Columns.Item(3).ColumnWidth = 109
$ws1.Columns.Item(5).ColumnWidth = 6.5

# --- Test Case Steps sheet: used range grows to include column D ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(9, 4).Value = "x"
$ws2.Cells.Item(9, 4).ClearContents()
